$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Apply updated values from the Feb 11 2023 15:50 UTC GitHub Actions refresh.
# Cells hold plain text (inline strings), so force Text number format before
# assigning to avoid Excel auto-converting numeric-looking / percent-looking
# strings into real numbers.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '308.96'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '0.82%'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '41.24'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '3.60%'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.130'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '1.69%'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.07644'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '0.03%'
$ws.Range('B6').NumberFormat = '@'
$ws.Range('B6').Value = 'GateToken'
$ws.Range('C6').NumberFormat = '@'
$ws.Range('C6').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '4.277'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '0.66%'
$ws.Range('B7').NumberFormat = '@'
$ws.Range('B7').Value = 'FTXToken'
$ws.Range('C7').NumberFormat = '@'
$ws.Range('C7').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.619'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '1.31%'
$ws.Range('B8').NumberFormat = '@'
$ws.Range('B8').Value = 'BTSEToken'
$ws.Range('C8').NumberFormat = '@'
$ws.Range('C8').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '2.473'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '2.23%'
$ws.Range('B9').NumberFormat = '@'
$ws.Range('B9').Value = 'MXToken'
$ws.Range('C9').NumberFormat = '@'
$ws.Range('C9').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.9087'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '-0.15%'
$ws.Range('B10').NumberFormat = '@'
$ws.Range('B10').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C10').NumberFormat = '@'
$ws.Range('C10').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.1177'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '13.47%'
$ws.Range('B11').NumberFormat = '@'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').NumberFormat = '@'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.1802'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '2.15%'
$ws.Range('B12').NumberFormat = '@'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').NumberFormat = '@'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.09186'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '-1.52%'
$ws.Range('B13').NumberFormat = '@'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').NumberFormat = '@'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.04261'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '-4.24%'
$ws.Range('B14').NumberFormat = '@'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').NumberFormat = '@'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.1043'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '-1.18%'
$ws.Range('B15').NumberFormat = '@'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').NumberFormat = '@'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.001257'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '-0.47%'
$ws.Range('B16').NumberFormat = '@'
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').NumberFormat = '@'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.005785'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '-1.18%'
$ws.Range('B17').NumberFormat = '@'
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').NumberFormat = '@'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.356'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '0.01%'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '0.47%'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.916'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '-0.40%'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.1375'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '1.92%'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.04028'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '-3.12%'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.001277'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '5.62%'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.004079'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '-0.19%'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '-2.55%'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0003747'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02433'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '-1.15%'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05256'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '2.07%'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.007807'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '-2.01%'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.1302'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '-0.22%'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.006790'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '-4.44%'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.001951'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '-0.26%'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.007543'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '-3.20%'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '0.71%'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '6.74%'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.00000000751'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '-0.27%'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.08095'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '1,674.80%'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.003003'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '-0.28%'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.00002102'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '-0.27%'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0002002'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '-0.27%'
